# Porcentajes_Integrantes.docx edit script
# 1) Turn the leading blank paragraph into a centered, bold, 14pt (sz=28)
#    title "PORCENTAJES OBTENIDOS DURANTE LA PRIMERA ENTREGA", with the
#    document's _GoBack bookmark sitting between the two title runs, and
#    add two more centered/bold/14pt blank paragraphs right after it.
# 2) Give the table an explicit preferred width and widen its two columns.
# 3) Clean up the "Shaid Bojorquez" cell: merge the "B" / "ojorquez" runs
#    (which used to be split apart by the old _GoBack bookmark) into a
#    single "Bojorquez" run, now that the bookmark lives in the title.

$d = $word.ActiveDocument

# --- 1) Title paragraph + two following blank (centered/bold/14pt) paragraphs ---
$p1 = $d.Paragraphs.Item(1)
$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>PORCENTAJES OBTENIDOS</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> DURANTE LA PRIMERA ENTREGA</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
[void]$p1.Range.InsertXML($titleXml)

# --- 2) Table width + column widths (values are OOXML dxa / 20 = points) ---
$t = $d.Tables.Item(1)
$t.PreferredWidthType = 3
$t.PreferredWidth = 9067 / 20
$t.Columns.Item(1).Width = 5387 / 20
$t.Columns.Item(2).Width = 3680 / 20

# --- 3) Merge the split "B" + "ojorquez" runs (drops the old bookmark there) ---
$rng = $d.Content
$found = $rng.Find.Execute("Bojorquez", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
  # direct same-text reassignment is a no-op, so stage through a placeholder
  $rng.Text = "__TEMP_NAME__"
  $rng2 = $d.Content
  $found2 = $rng2.Find.Execute("__TEMP_NAME__", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
  if ($found2) {
    $rng2.Text = "Bojorquez"
  }
}

Write-Output "edit complete"
